# Update strategy to exit at support/resistance levels
# Adjusts Exit price, Exit Reason, and PnL for trades that exited at a
# support or resistance level instead of at the closing price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 -> Resistance Level
$ws.Range("C3").Value = 6696.18017578125
$ws.Range("E3").Value = "Resistance Level"
$ws.Range("F3").Value = 30.80029296875

# Row 4 -> Resistance Level
$ws.Range("C4").Value = 6808.9501953125
$ws.Range("E4").Value = "Resistance Level"
$ws.Range("F4").Value = 22.1103515625

# Row 5 -> Resistance Level
$ws.Range("C5").Value = 6827.31005859375
$ws.Range("E5").Value = "Resistance Level"
$ws.Range("F5").Value = 10.240234375

# Row 16 -> Support Level
$ws.Range("C16").Value = 6840.60986328125
$ws.Range("E16").Value = "Support Level"
$ws.Range("F16").Value = 13.5302734375

# Row 20 -> Support Level
$ws.Range("C20").Value = 6771.72021484375
$ws.Range("E20").Value = "Support Level"
$ws.Range("F20").Value = 33.27001953125

# Row 22 -> Resistance Level
$ws.Range("C22").Value = 6912.6201171875
$ws.Range("E22").Value = "Resistance Level"
$ws.Range("F22").Value = 11.64013671875

# Row 23 -> Support Level
$ws.Range("C23").Value = 6895.2998046875
$ws.Range("E23").Value = "Support Level"
$ws.Range("F23").Value = 4.85009765625

# Row 26 -> Support Level
$ws.Range("C26").Value = 6960.81005859375
$ws.Range("E26").Value = "Support Level"
$ws.Range("F26").Value = 11.2900390625
